$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.167.39'
$ws.Range("E2").Value = '  -2.50%  '

$ws.Range("D3").Value = '1.722.96'
$ws.Range("E3").Value = '  -2.52%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.63'
$ws.Range("E5").Value = '  -4.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4634'
$ws.Range("E7").Value = '  +3.46%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3446'
$ws.Range("E8").Value = '  -3.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.63'
$ws.Range("E9").Value = '  +1.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07306'
$ws.Range("E10").Value = '  -1.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.054'
$ws.Range("E11").Value = '  -3.78%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9990'
$ws.Range("E12").Value = '  -0.18%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.89'
$ws.Range("E13").Value = '  -4.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.866'
$ws.Range("E14").Value = '  -2.79%  '

$ws.Range("D15").Value = '1.717.18'
$ws.Range("E15").Value = '  -3.23%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.910'
$ws.Range("E16").Value = '  -4.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.71'
$ws.Range("E17").Value = '  -3.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001044'
$ws.Range("E18").Value = '  -1.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06323'
$ws.Range("E19").Value = '  -1.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.52'
$ws.Range("E21").Value = '  -3.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.656'
$ws.Range("E22").Value = '  -2.59%  '

$ws.Range("D23").Value = '27.207.52'
$ws.Range("E23").Value = '  -2.52%  '

$ws.Range("E24").Value = '  -3.89%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.137'
$ws.Range("E25").Value = '  +1.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.38'
$ws.Range("E26").Value = '  -3.82%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.49'
$ws.Range("E27").Value = '  -3.62%  '

$ws.Range("D28").Value = '1.919.32'
$ws.Range("E28").Value = '  -2.95%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.156'
$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.33'
$ws.Range("E30").Value = '  -4.80%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.036'
$ws.Range("E31").Value = '  -5.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09114'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.582'
$ws.Range("E33").Value = '  -1.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.363'
$ws.Range("E34").Value = '  -3.94%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02204'
$ws.Range("E35").Value = '  -4.00%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05865'
$ws.Range("E36").Value = '  -3.82%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.16'
$ws.Range("E37").Value = '  -5.61%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.1995'
$ws.Range("E38").Value = '  -4.73%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.767'
$ws.Range("E39").Value = '  -3.95%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5984'
$ws.Range("E40").Value = '  -5.36%  '

$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.399'
$ws.Range("E41").Value = '  +0.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.133'
$ws.Range("E42").Value = '  -4.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.527'
$ws.Range("E43").Value = '  -4.96%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.67'
$ws.Range("E44").Value = '  -4.27%  '

$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.620'
$ws.Range("E45").Value = '  -3.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5631'
$ws.Range("E46").Value = '  -4.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '119.22'
$ws.Range("E47").Value = '  -2.61%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.871'
$ws.Range("E48").Value = '  -4.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06665'
$ws.Range("E49").Value = '  -3.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.090'
$ws.Range("E50").Value = '  -4.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9989'
$ws.Range("E51").Value = '  -0.12%  '
